$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(94, 8).Value = 4942
$ws.Cells.Item(94, 9).Value = 1570
$ws.Cells.Item(94, 10).Value = 10000
$ws.Cells.Item(94, 11).Value = 1570
$ws.Cells.Item(94, 12).Value = 10000
$ws.Cells.Item(94, 13).Value = -1119
$ws.Cells.Item(94, 14).Value = -10902

$ws.Cells.Item(100, 8).Value = 2612.5
$ws.Cells.Item(100, 9).Value = 2475
$ws.Cells.Item(100, 10).Value = 2750
$ws.Cells.Item(100, 11).Value = 2475
$ws.Cells.Item(100, 12).Value = 2750
$ws.Cells.Item(100, 13).Value = -1934
$ws.Cells.Item(100, 14).Value = -3832

$ws.Cells.Item(113, 8).Value = 3003.5557
$ws.Cells.Item(113, 9).Value = 2998.6667
$ws.Cells.Item(113, 10).Value = 3006
$ws.Cells.Item(113, 11).Value = 2998.6667
$ws.Cells.Item(113, 12).Value = 3006
$ws.Cells.Item(113, 13).Value = 255.3332999999998
$ws.Cells.Item(113, 14).Value = -9514

$ws.Cells.Item(116, 8).Value = 8691.6
$ws.Cells.Item(116, 9).Value = 13200.444
$ws.Cells.Item(116, 10).Value = 1928.3334
$ws.Cells.Item(116, 11).Value = 13200.444
$ws.Cells.Item(116, 12).Value = 1928.3334
$ws.Cells.Item(116, 13).Value = -9758.444
$ws.Cells.Item(116, 14).Value = -8812.3334

$ws.Cells.Item(132, 8).Value = 2644.3125
$ws.Cells.Item(132, 9).Value = 2540.2932
$ws.Cells.Item(132, 10).Value = 3649.8333
$ws.Cells.Item(132, 11).Value = 7620.8796
$ws.Cells.Item(132, 12).Value = 10949.4999
$ws.Cells.Item(132, 13).Value = -5090.8796
$ws.Cells.Item(132, 14).Value = -16009.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 3367
$ws.Cells.Item(45, 9).Value = 3389.3333
$ws.Cells.Item(45, 10).Value = 3300
$ws.Cells.Item(45, 11).Value = 3389.3333
$ws.Cells.Item(45, 12).Value = 3300
$ws.Cells.Item(45, 13).Value = -3012.3333
$ws.Cells.Item(45, 14).Value = -4054

$ws.Cells.Item(61, 8).Value = 2856.7576
$ws.Cells.Item(61, 9).Value = 2133
$ws.Cells.Item(61, 10).Value = 4304.273
$ws.Cells.Item(61, 11).Value = 2133
$ws.Cells.Item(61, 12).Value = 4304.273
$ws.Cells.Item(61, 13).Value = -1921
$ws.Cells.Item(61, 14).Value = -4728.273

$ws.Cells.Item(74, 8).Value = 1847.1111
$ws.Cells.Item(74, 9).Value = 1453
$ws.Cells.Item(74, 10).Value = 5000
$ws.Cells.Item(74, 11).Value = 1453
$ws.Cells.Item(74, 12).Value = 5000
$ws.Cells.Item(74, 13).Value = -579
$ws.Cells.Item(74, 14).Value = -6748

$ws.Cells.Item(77, 8).Value = 1847.1111
$ws.Cells.Item(77, 9).Value = 1453
$ws.Cells.Item(77, 10).Value = 5000
$ws.Cells.Item(77, 11).Value = 7265
$ws.Cells.Item(77, 12).Value = 25000
$ws.Cells.Item(77, 13).Value = -2897
$ws.Cells.Item(77, 14).Value = -33736

$ws.Cells.Item(132, 8).Value = 2509.016
$ws.Cells.Item(132, 9).Value = 1823.3489
$ws.Cells.Item(132, 10).Value = 4060.7896
$ws.Cells.Item(132, 11).Value = 5470.0467
$ws.Cells.Item(132, 12).Value = 12182.3688
$ws.Cells.Item(132, 13).Value = -2940.0467
$ws.Cells.Item(132, 14).Value = -17242.3688

$ws.Cells.Item(136, 8).Value = 2856.7576
$ws.Cells.Item(136, 9).Value = 2133
$ws.Cells.Item(136, 10).Value = 4304.273
$ws.Cells.Item(136, 11).Value = 6399
$ws.Cells.Item(136, 12).Value = 12912.819
$ws.Cells.Item(136, 13).Value = -3849
$ws.Cells.Item(136, 14).Value = -18012.819

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2586.4167
$ws.Cells.Item(86, 9).Value = 2730.7778
$ws.Cells.Item(86, 10).Value = 2153.3333
$ws.Cells.Item(86, 11).Value = 2730.7778
$ws.Cells.Item(86, 12).Value = 2153.3333
$ws.Cells.Item(86, 13).Value = -1607.7778
$ws.Cells.Item(86, 14).Value = -4399.3333

$ws.Cells.Item(89, 8).Value = 2586.4167
$ws.Cells.Item(89, 9).Value = 2730.7778
$ws.Cells.Item(89, 10).Value = 2153.3333
$ws.Cells.Item(89, 11).Value = 13653.889
$ws.Cells.Item(89, 12).Value = 10766.6665
$ws.Cells.Item(89, 13).Value = -8037.888999999999
$ws.Cells.Item(89, 14).Value = -21998.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(75, 8).Value = 80500
$ws.Cells.Item(75, 9).Value = 0
$ws.Cells.Item(75, 10).Value = 80500
$ws.Cells.Item(75, 11).Value = 0
$ws.Cells.Item(75, 12).Value = 80500
$ws.Cells.Item(75, 14).Value = -82496

$ws.Cells.Item(78, 8).Value = 80500
$ws.Cells.Item(78, 9).Value = 0
$ws.Cells.Item(78, 10).Value = 80500
$ws.Cells.Item(78, 11).Value = 0
$ws.Cells.Item(78, 12).Value = 241500
$ws.Cells.Item(78, 14).Value = -251484

$ws.Cells.Item(81, 8).Value = 0
$ws.Cells.Item(81, 9).Value = 0
$ws.Cells.Item(81, 10).Value = 0
$ws.Cells.Item(81, 11).Value = 0
$ws.Cells.Item(81, 12).Value = 0
$ws.Cells.Item(81, 14).Value = ""

$ws.Cells.Item(84, 8).Value = 0
$ws.Cells.Item(84, 9).Value = 0
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 11).Value = 0
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 14).Value = ""

$ws.Cells.Item(99, 8).Value = 2033.3704
$ws.Cells.Item(99, 9).Value = 1625.5
$ws.Cells.Item(99, 10).Value = 2104.3044
$ws.Cells.Item(99, 11).Value = 1625.5
$ws.Cells.Item(99, 12).Value = 2104.3044
$ws.Cells.Item(99, 13).Value = -127.5
$ws.Cells.Item(99, 14).Value = -5100.3044

$ws.Cells.Item(126, 8).Value = 2033.3704
$ws.Cells.Item(126, 9).Value = 1625.5
$ws.Cells.Item(126, 10).Value = 2104.3044
$ws.Cells.Item(126, 11).Value = 4876.5
$ws.Cells.Item(126, 12).Value = 6312.9132
$ws.Cells.Item(126, 13).Value = -2406.5
$ws.Cells.Item(126, 14).Value = -11252.9132

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 486.5098
$ws.Cells.Item(113, 9).Value = 487.14285
$ws.Cells.Item(113, 10).Value = 486.27026
$ws.Cells.Item(113, 11).Value = 1461.42855
$ws.Cells.Item(113, 12).Value = 1458.81078
$ws.Cells.Item(113, 13).Value = 708.5714499999999
$ws.Cells.Item(113, 14).Value = -5798.81078

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 9466.200000000001
$ws.Cells.Item(70, 9).Value = 10908.267
$ws.Cells.Item(70, 10).Value = 5140
$ws.Cells.Item(70, 11).Value = 10908.267
$ws.Cells.Item(70, 12).Value = 5140
$ws.Cells.Item(70, 13).Value = -10638.267
$ws.Cells.Item(70, 14).Value = -5680

$ws.Cells.Item(73, 8).Value = 9466.200000000001
$ws.Cells.Item(73, 9).Value = 10908.267
$ws.Cells.Item(73, 10).Value = 5140
$ws.Cells.Item(73, 11).Value = 10908.267
$ws.Cells.Item(73, 12).Value = 5140
$ws.Cells.Item(73, 13).Value = -9972.267
$ws.Cells.Item(73, 14).Value = -7012

$ws.Cells.Item(94, 8).Value = 41848
$ws.Cells.Item(94, 9).Value = 0
$ws.Cells.Item(94, 10).Value = 41848
$ws.Cells.Item(94, 11).Value = 0
$ws.Cells.Item(94, 12).Value = 41848
$ws.Cells.Item(94, 14).Value = -43200

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5493.615
$ws.Cells.Item(7, 9).Value = 4712.4443
$ws.Cells.Item(7, 10).Value = 7251.25
$ws.Cells.Item(7, 11).Value = 4712.4443
$ws.Cells.Item(7, 12).Value = 7251.25
$ws.Cells.Item(7, 13).Value = -4600.4443
$ws.Cells.Item(7, 14).Value = -7475.25

$ws.Cells.Item(68, 8).Value = 1533.8077
$ws.Cells.Item(68, 9).Value = 1475.1904
$ws.Cells.Item(68, 10).Value = 1780
$ws.Cells.Item(68, 11).Value = 1475.1904
$ws.Cells.Item(68, 12).Value = 1780
$ws.Cells.Item(68, 13).Value = -726.1904
$ws.Cells.Item(68, 14).Value = -3278

$ws.Cells.Item(71, 8).Value = 1533.8077
$ws.Cells.Item(71, 9).Value = 1475.1904
$ws.Cells.Item(71, 10).Value = 1780
$ws.Cells.Item(71, 11).Value = 7375.951999999999
$ws.Cells.Item(71, 12).Value = 8900
$ws.Cells.Item(71, 13).Value = -3631.951999999999
$ws.Cells.Item(71, 14).Value = -16388

$ws.Cells.Item(93, 8).Value = 4606.148
$ws.Cells.Item(93, 9).Value = 6518.294
$ws.Cells.Item(93, 10).Value = 1355.5
$ws.Cells.Item(93, 11).Value = 6518.294
$ws.Cells.Item(93, 12).Value = 1355.5
$ws.Cells.Item(93, 13).Value = -5270.294
$ws.Cells.Item(93, 14).Value = -3851.5

$ws.Cells.Item(100, 8).Value = 2720.5
$ws.Cells.Item(100, 9).Value = 2637.7144
$ws.Cells.Item(100, 10).Value = 3300
$ws.Cells.Item(100, 11).Value = 2637.7144
$ws.Cells.Item(100, 12).Value = 3300
$ws.Cells.Item(100, 13).Value = -2096.7144
$ws.Cells.Item(100, 14).Value = -4382

$ws.Cells.Item(126, 8).Value = 5493.615
$ws.Cells.Item(126, 9).Value = 4712.4443
$ws.Cells.Item(126, 10).Value = 7251.25
$ws.Cells.Item(126, 11).Value = 14137.3329
$ws.Cells.Item(126, 12).Value = 21753.75
$ws.Cells.Item(126, 13).Value = -11667.3329
$ws.Cells.Item(126, 14).Value = -26693.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 53429
$ws.Cells.Item(46, 9).Value = 0
$ws.Cells.Item(46, 10).Value = 53429
$ws.Cells.Item(46, 11).Value = 0
$ws.Cells.Item(46, 12).Value = 53429
$ws.Cells.Item(46, 14).Value = -53891

$ws.Cells.Item(81, 8).Value = 3958
$ws.Cells.Item(81, 9).Value = 4441.1113
$ws.Cells.Item(81, 10).Value = 3233.3333
$ws.Cells.Item(81, 11).Value = 8882.222599999999
$ws.Cells.Item(81, 12).Value = 6466.6666
$ws.Cells.Item(81, 13).Value = -7821.222599999999
$ws.Cells.Item(81, 14).Value = -8588.6666

$ws.Cells.Item(84, 8).Value = 3958
$ws.Cells.Item(84, 9).Value = 4441.1113
$ws.Cells.Item(84, 10).Value = 3233.3333
$ws.Cells.Item(84, 11).Value = 44411.113
$ws.Cells.Item(84, 12).Value = 32333.333
$ws.Cells.Item(84, 13).Value = -39107.113
$ws.Cells.Item(84, 14).Value = -42941.333

$ws.Cells.Item(134, 8).Value = 53429
$ws.Cells.Item(134, 9).Value = 0
$ws.Cells.Item(134, 10).Value = 53429
$ws.Cells.Item(134, 11).Value = 0
$ws.Cells.Item(134, 12).Value = 160287
$ws.Cells.Item(134, 14).Value = -165357

$ws.Cells.Item(135, 8).Value = 111592.5
$ws.Cells.Item(135, 9).Value = 0
$ws.Cells.Item(135, 10).Value = 111592.5
$ws.Cells.Item(135, 11).Value = 0
$ws.Cells.Item(135, 12).Value = 111592.5
$ws.Cells.Item(135, 14).Value = -121732.5
